$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns we touch so Excel
# does not auto-convert numeric-looking strings (e.g. "2.80", "11.50")
# into numbers and strip formatting - source cells are inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "65.042.93"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3
$ws.Range("D3").Value = "3.363.17"
$ws.Range("E3").Value = "  +1.43%  "

# Row 4
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "181.25"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "536.51"
$ws.Range("E6").Value = "  +1.19%  "

# Row 7
$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  -1.23%  "

# Row 8
$ws.Range("D8").Value = "3.348.32"
$ws.Range("E8").Value = "  +1.20%  "

# Row 9
$ws.Range("E9").Value = "  -0.23%  "

# Row 10
$ws.Range("D10").Value = "0.617"
$ws.Range("E10").Value = "  +0.62%  "

# Row 11
$ws.Range("D11").Value = "55.28"
$ws.Range("E11").Value = "  -7.05%  "

# Row 12
$ws.Range("D12").Value = "0.139"
$ws.Range("E12").Value = "  +4.07%  "

# Row 13
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  +1.72%  "

# Row 14
$ws.Range("D14").Value = "9.11"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").Value = "3.850.90"
$ws.Range("E15").Value = "  -0.05%  "

# Row 16
$ws.Range("D16").Value = "0.119"
$ws.Range("E16").Value = "  +1.25%  "

# Row 17
$ws.Range("D17").Value = "3.326.57"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").Value = "64.985.19"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19
$ws.Range("D19").Value = "17.85"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("D20").Value = "11.34"
$ws.Range("E20").Value = "  +1.21%  "

# Row 21
$ws.Range("D21").Value = "0.976"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$ws.Range("D22").Value = "386.01"
$ws.Range("E22").Value = "  +2.64%  "

# Row 23
$ws.Range("D23").Value = "4.21"
$ws.Range("E23").Value = "  +7.16%  "

# Row 24
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.65"
$ws.Range("E24").Value = "  +4.41%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "82.63"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("D26").Value = "3.76"
$ws.Range("E26").Value = "  -1.89%  "

# Row 27
$ws.Range("D27").Value = "6.11"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("D28").Value = "2.80"
$ws.Range("E28").Value = "  +3.72%  "

# Row 29
$ws.Range("D29").Value = "11.50"
$ws.Range("E29").Value = "  -0.84%  "

# Row 30
$ws.Range("D30").Value = "8.39"
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
$ws.Range("D31").Value = "29.29"
$ws.Range("E31").Value = "  +0.50%  "

# Row 32
$ws.Range("D32").Value = "648.36"
$ws.Range("E32").Value = "  -1.04%  "

# Row 33
$ws.Range("D33").Value = "6.78"
$ws.Range("E33").Value = "  +0.93%  "

# Row 34
$ws.Range("D34").Value = "11.36"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  +1.21%  "

# Row 36
$ws.Range("D36").Value = "57.78"
$ws.Range("E36").Value = "  -2.75%  "

# Row 37
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.26%  "

# Row 38
$ws.Range("D38").Value = "37.35"
$ws.Range("E38").Value = "  +1.31%  "

# Row 39
$ws.Range("D39").Value = "0.390"
$ws.Range("E39").Value = "  -1.26%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  +10.60%  "

# Row 41
$ws.Range("D41").Value = "0.995"
$ws.Range("E41").Value = "  -0.33%  "

# Row 42
$ws.Range("D42").Value = "3.27"
$ws.Range("E42").Value = "  +14.68%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +8.25%  "

# Row 44
$ws.Range("D44").Value = "0.128"
$ws.Range("E44").Value = "  +0.97%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.004.42"
$ws.Range("E45").Value = "  +3.72%  "

# Row 46
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  +1.29%  "

# Row 47
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.20"
$ws.Range("E48").Value = "  +3.72%  "

# Row 49
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.69"
$ws.Range("E49").Value = "  +1.50%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.127"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "138.92"
$ws.Range("E51").Value = "  +2.25%  "
